# sync db when saving contract
# Update the generated contract fields to reflect the latest values
# pulled from the database before the document is saved.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# تعريف السيارة (vehicle identification)
Replace-Text "renault   " "vp   "
Replace-Text "r3   " "ppppppppp   "
Replace-Text "voila " "ppppppp "
Replace-Text "333111212  " "999999999  "
Replace-Text "455555 " "999999 "

# مؤشر السير - must run before the "مبلغ الضمان" (0) replacement below,
# since "4500   " itself contains the substring "0   ".
Replace-Text "4500   " "999999   "

# السعر / مبلغ الضمان / جواز السفر
Replace-Text "150  " "999  "
Replace-Text "0   " "9999999   "
Replace-Text "1144   " "2222222   "

# اسم المستفيد / تاريخ الميلاد / الاسم
Replace-Text "bouaf   " "hamid hamid   "
Replace-Text "14/12/2023   " "15/12/2023   "
Replace-Text "husein" "aaaaaa"

# العنوان / رقم الهاتف
Replace-Text "aazaze   " "aaaaaaaaaa   "
Replace-Text "055555   " "000110   "

# رخصة السياقة: رقمها، تاريخ ومكان الإصدار
Replace-Text "111111" "111101010"
Replace-Text "08/12/2023" "14/12/2023"
Replace-Text "bordj  " "aaaaaaa  "

# التسليم: تاريخ، ساعة، المدة، تاريخ العودة
Replace-Text "03/12/2023" "09/12/2023"
Replace-Text "14:06   " "20:38   "
Replace-Text "لمدة:   2" "لمدة:   9"
Replace-Text "05/12/2023" "18/12/2023"
